$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: use existing rows 78:79 (shaded "WT-less" template) and 80:81
#     (non-shaded template) as formatting donors for the brand new rows
#     that will be appended at the bottom of the table (82:83 and 84:85),
#     copying them *before* their source content is overwritten below.
$ws.Range("A78:D79").Copy($ws.Range("A82:D83"))
$ws.Range("A80:D81").Copy($ws.Range("A84:D85"))

# --- Step 2: fix a stray mistake in row 77 (experiment id typo)
$ws.Cells.Item(77, 4).Value = "NLA_096"

# --- Step 3: row 78/79 used to describe the "(no TG)" condition for
#     experiment NLA_915; they actually describe the new "Tom - Age (WT)"
#     condition, so correct their text + shift the experiment id that had
#     been misassigned to NLA_916.
$ws.Cells.Item(78, 2).Value = "Tom – Age (WT)"
$ws.Cells.Item(78, 3).Value = "Genotype C22 – WT"
$ws.Cells.Item(78, 4).Value = "NLA_095"

$ws.Cells.Item(79, 2).Value = "Tom – Age (WT)"
$ws.Cells.Item(79, 3).Value = "Genotype C3 – WT"
$ws.Cells.Item(79, 4).Value = "NLA_097"

# --- Step 4: row 80/81 actually hold the "Tom (no TG)" experiment
#     (NLA_916), not "Tom - Age (no TG)"; fix the title text and drop the
#     shading so they match the other non-shaded "base" rows, as well as
#     bump the trailing experiment ids.
$ws.Cells.Item(80, 2).Value = "Tom (no TG)"
$ws.Cells.Item(80, 4).Value = "NLA_099"

$ws.Cells.Item(81, 2).Value = "Tom (no TG)"
$ws.Cells.Item(81, 4).Value = "NLA_100"

# --- Step 5: populate the newly appended rows.
# Rows 82:83 -> NLA_917 "Tom - Age (no TG)" (shaded, like the old 78:79)
$ws.Cells.Item(82, 1).Value = "NLA_917"
$ws.Cells.Item(82, 2).Value = "Tom – Age (no TG)"
$ws.Cells.Item(82, 3).Value = "Genotype C22 (no TG)"
$ws.Cells.Item(82, 4).Value = "NLA_101"

$ws.Cells.Item(83, 1).Value = "NLA_917"
$ws.Cells.Item(83, 2).Value = "Tom – Age (no TG)"
$ws.Cells.Item(83, 3).Value = "Genotype C3 (no TG)"
$ws.Cells.Item(83, 4).Value = "NLA_103"

# Rows 84:85 -> NLA_918 "Tom - Age (no TG) (WT)" (non-shaded, like 80:81)
$ws.Cells.Item(84, 1).Value = "NLA_918"
$ws.Cells.Item(84, 2).Value = "Tom – Age (no TG) (WT)"
$ws.Cells.Item(84, 3).Value = "Genotype C22 – WT (no TG)"
$ws.Cells.Item(84, 4).Value = "NLA_102"

$ws.Cells.Item(85, 1).Value = "NLA_918"
$ws.Cells.Item(85, 2).Value = "Tom – Age (no TG) (WT)"
$ws.Cells.Item(85, 3).Value = "Genotype C3 – WT (no TG)"
$ws.Cells.Item(85, 4).Value = "NLA_104"

# --- Step 6: view bookkeeping (zoom + scroll position + active selection)
$excel.ActiveWindow.Zoom = 75
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 47
$ws.Range("D86").Select()
